## report_2022_monthly_rates.xlsx — "wire transfer" -> "currency conversion"
## rewording + the knock-on numeric corrections that came with it.
##
## Changes applied (mirrors the target OOXML diff):
##   1. "Fees" sheet: both cells that said "Wire transfer" (B6 and B11) are
##      reworded to "Currency conversion or wire transfer" (the underlying
##      shared-string text change).
##   2. "Fees" sheet: column B is widened so the new, longer label still
##      fits (bestFit-style width).
##   3. "Foreign Currencies" sheet, row 7: the USD amount converted (B7) and
##      the resulting EUR gain/loss (G7) are corrected.
##   4. "ELSTER - Summary" sheet, row 7 (the currency gain/loss summary
##      line that is fed by the "Foreign Currencies" sheet): C7 is updated
##      to match the corrected totals.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1 & 2. "Fees" sheet — reword the "Wire transfer" label and widen col B
# ---------------------------------------------------------------------
$fees = $wb.Worksheets.Item("Fees")

$newLabel = "Currency conversion or wire transfer"
$fees.Range("B6").Value2 = $newLabel
$fees.Range("B11").Value2 = $newLabel

# Widen column B so the longer text still "best fits" (was 12.85546875).
$fees.Columns.Item(2).ColumnWidth = 32.3

# ---------------------------------------------------------------------
# 3. "Foreign Currencies" sheet, row 7 — corrected conversion amount and
#    EUR gain/loss
# ---------------------------------------------------------------------
$foreignCurrencies = $wb.Worksheets.Item("Foreign Currencies")
$foreignCurrencies.Range("B7").Value2 = 155
$foreignCurrencies.Range("G7").Value2 = -10.13

# ---------------------------------------------------------------------
# 4. "ELSTER - Summary" sheet, row 7 — updated currency gain/loss total
# ---------------------------------------------------------------------
$elsterSummary = $wb.Worksheets.Item("ELSTER - Summary")
$elsterSummary.Range("C7").Value2 = 16.86
